$wb = $excel.ActiveWorkbook

# --- "Expected Out" sheet: bump a couple of amounts (B1 SUM formula recalcs automatically) ---
$expectedOut = $wb.Worksheets.Item("Expected Out")
$expectedOut.Range("B9").Value = 1345.76
$expectedOut.Range("B11").Value = 425.62

# --- "Budget Out" sheet: bump an amount and lengthen a description string ---
$budgetOut = $wb.Worksheets.Item("Budget Out")
$budgetOut.Range("C9").Value = 87.82
$budgetOut.Range("F9").Value = "Description007zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# --- "TestRecord" sheet: move a transaction's date forward a few days and bump its amount ---
$testRecord = $wb.Worksheets.Item("TestRecord")
$testRecord.Range("A10").Value = 43261
$testRecord.Range("B10").Value = 118.74
$testRecord.Range("E10").Value = "some test textzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

$wb.Save()
